$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (row 11): correct-answer marking value 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row (row 12): total correct marks 24 -> 40 (8 correct * 5 marks)
$ws.Range("B12").Value = 40

# Update the corr/total marks fraction text 23/84 -> 40/140
$ws.Range("E12").Value = "40/140"
